# chore: update Sheets via scheduled runner
# Refresh cached market/profit figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# for a handful of leves across the job sheets, matching a scheduled market-data sync.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 225.42105
$ws.Range("I53").Value = 130.14285
$ws.Range("J53").Value = 281
$ws.Range("K53").Value = 130.14285
$ws.Range("L53").Value = 281
$ws.Range("M53").Value = 506.85715
$ws.Range("N53").Value = -1555
$ws.Range("H125").Value = 3104
$ws.Range("I125").Value = 700
$ws.Range("J125").Value = 3638.2222
$ws.Range("K125").Value = 6300
$ws.Range("L125").Value = 32743.9998
$ws.Range("M125").Value = -3840
$ws.Range("N125").Value = -37663.99980000001
$ws.Range("H141").Value = 26941232
$ws.Range("I141").Value = 38650910
$ws.Range("J141").Value = 8970.5
$ws.Range("K141").Value = 115952730
$ws.Range("L141").Value = 26911.5
$ws.Range("M141").Value = -115947550
$ws.Range("N141").Value = -37271.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1759.93
$ws.Range("I32").Value = 1759.93
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1759.93
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1472.93
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 1472.2174
$ws.Range("I61").Value = 1182.2565
$ws.Range("J61").Value = 3087.7144
$ws.Range("K61").Value = 1182.2565
$ws.Range("L61").Value = 3087.7144
$ws.Range("M61").Value = -970.2565
$ws.Range("N61").Value = -3511.7144
$ws.Range("H63").Value = 2310.076
$ws.Range("I63").Value = 2293.5483
$ws.Range("J63").Value = 2370.353
$ws.Range("K63").Value = 2293.5483
$ws.Range("L63").Value = 2370.353
$ws.Range("M63").Value = -1607.5483
$ws.Range("N63").Value = -3742.353
$ws.Range("H66").Value = 2310.076
$ws.Range("I66").Value = 2293.5483
$ws.Range("J66").Value = 2370.353
$ws.Range("K66").Value = 11467.7415
$ws.Range("L66").Value = 11851.765
$ws.Range("M66").Value = -8035.7415
$ws.Range("N66").Value = -18715.765
$ws.Range("H102").Value = 1740
$ws.Range("I102").Value = 1588.8889
$ws.Range("J102").Value = 3100
$ws.Range("K102").Value = 1588.8889
$ws.Range("L102").Value = 3100
$ws.Range("M102").Value = 33.11110000000008
$ws.Range("N102").Value = -6344
$ws.Range("H136").Value = 1472.2174
$ws.Range("I136").Value = 1182.2565
$ws.Range("J136").Value = 3087.7144
$ws.Range("K136").Value = 3546.7695
$ws.Range("L136").Value = 9263.143199999999
$ws.Range("M136").Value = -996.7694999999999
$ws.Range("N136").Value = -14363.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1931.875
$ws.Range("I86").Value = 1679.875
$ws.Range("J86").Value = 2435.875
$ws.Range("K86").Value = 1679.875
$ws.Range("L86").Value = 2435.875
$ws.Range("M86").Value = -556.875
$ws.Range("N86").Value = -4681.875
$ws.Range("H89").Value = 1931.875
$ws.Range("I89").Value = 1679.875
$ws.Range("J89").Value = 2435.875
$ws.Range("K89").Value = 8399.375
$ws.Range("L89").Value = 12179.375
$ws.Range("M89").Value = -2783.375
$ws.Range("N89").Value = -23411.375
$ws.Range("H99").Value = 1810.3334
$ws.Range("I99").Value = 1954.1666
$ws.Range("J99").Value = 1666.5
$ws.Range("K99").Value = 1954.1666
$ws.Range("L99").Value = 1666.5
$ws.Range("M99").Value = -456.1666
$ws.Range("N99").Value = -4662.5
$ws.Range("H105").Value = 758931.75
$ws.Range("I105").Value = 948152.2
$ws.Range("J105").Value = 2050
$ws.Range("K105").Value = 948152.2
$ws.Range("L105").Value = 2050
$ws.Range("M105").Value = -946405.2
$ws.Range("N105").Value = -5544
$ws.Range("H134").Value = 4465822
$ws.Range("I134").Value = 6580186.5
$ws.Range("J134").Value = 2163.2222
$ws.Range("K134").Value = 19740559.5
$ws.Range("L134").Value = 6489.6666
$ws.Range("M134").Value = -19738024.5
$ws.Range("N134").Value = -11559.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1447.9474
$ws.Range("I107").Value = 585.46155
$ws.Range("J107").Value = 3316.6667
$ws.Range("K107").Value = 585.46155
$ws.Range("L107").Value = 3316.6667
$ws.Range("M107").Value = 1334.53845
$ws.Range("N107").Value = -7156.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 534.63635
$ws.Range("I50").Value = 105
$ws.Range("J50").Value = 695.75
$ws.Range("K50").Value = 315
$ws.Range("L50").Value = 2087.25
$ws.Range("M50").Value = 166
$ws.Range("N50").Value = -3049.25
$ws.Range("H53").Value = 534.63635
$ws.Range("I53").Value = 105
$ws.Range("J53").Value = 695.75
$ws.Range("K53").Value = 315
$ws.Range("L53").Value = 2087.25
$ws.Range("M53").Value = 166
$ws.Range("N53").Value = -3049.25
$ws.Range("H136").Value = 785
$ws.Range("I136").Value = 585.7143
$ws.Range("K136").Value = 1757.1429
$ws.Range("M136").Value = 3342.8571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 736.5454999999999
$ws.Range("I31").Value = 736.5454999999999
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 736.5454999999999
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -444.5454999999999
$ws.Range("N31").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H37").Value = 736.5454999999999
$ws.Range("I37").Value = 736.5454999999999
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 736.5454999999999
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -459.5454999999999
$ws.Range("N37").ClearContents()
$ws.Range("H113").Value = 1935.0952
$ws.Range("I113").Value = 1929.3572
$ws.Range("J113").Value = 1946.5714
$ws.Range("K113").Value = 1929.3572
$ws.Range("L113").Value = 1946.5714
$ws.Range("M113").Value = 240.6428000000001
$ws.Range("N113").Value = -6286.5714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 17585.8
$ws.Range("J112").Value = 17585.8
$ws.Range("L112").Value = 17585.8
$ws.Range("N112").Value = -20539.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1831.0769
$ws.Range("I122").Value = 1511.5555
$ws.Range("J122").Value = 2550
$ws.Range("K122").Value = 4534.666499999999
$ws.Range("L122").Value = 7650
$ws.Range("M122").Value = -2084.666499999999
$ws.Range("N122").Value = -12550

